# "Generate Report for Archive" — refresh the localization-status report:
#   - flip the in-flight rows from "Ready for handoff" to "In Translation"
#   - re-fit the Status/locale columns now that the text is shorter

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Column widths shrink to fit the new, shorter status text.
$overview.Columns("E:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
